$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("A1").Value = [double]"0.8148776013935055"
$ws.Range("B1").Value = [double]"0.04508152606431248"
$ws.Range("A2").Value = [double]"0.1520929998521841"
$ws.Range("B2").Value = [double]"0.7779131012457334"
$ws.Range("C2").Value = [double]"0.7859027136591488"
$ws.Range("D2").Value = [double]"0.2600665302161255"
$ws.Range("A3").Value = [double]"1.099127094764975e-18"
$ws.Range("B3").Value = [double]"5.6108765779337e-18"
$ws.Range("C3").Value = [double]"5.626651342498006e-18"
$ws.Range("D3").Value = [double]"1.783871627779743e-18"
$ws.Range("A4").Value = [double]"3.592543686251221e-18"
$ws.Range("B4").Value = [double]"1.834310814168978e-17"
$ws.Range("C4").Value = [double]"1.840898925186547e-17"
$ws.Range("D4").Value = [double]"5.862326623254954e-18"
$ws.Range("A5").Value = [double]"0.03302938385314936"
$ws.Range("B5").Value = [double]"0.1770053950416967"
$ws.Range("C5").Value = [double]"0.2140972714396903"
$ws.Range("D5").Value = [double]"0.7399334668404353"

$ws = $wb.Worksheets.Item(2)
$ws.Range("A1").Value = [double]"0.8664976520894124"
$ws.Range("B1").Value = [double]"0.08736546971699768"
$ws.Range("A2").Value = [double]"0.1057632657938702"
$ws.Range("B2").Value = [double]"0.7170620091882586"
$ws.Range("C2").Value = [double]"0.7562781713913763"
$ws.Range("D2").Value = [double]"0.2633149498588686"
$ws.Range("A3").Value = [double]"1.018502906551284e-18"
$ws.Range("B3").Value = [double]"6.680678484047883e-18"
$ws.Range("C3").Value = [double]"6.127780568523731e-18"
$ws.Range("D3").Value = [double]"1.079730753712121e-18"
$ws.Range("A4").Value = [double]"2.050329005769312e-19"
$ws.Range("B4").Value = [double]"1.344963883311013e-18"
$ws.Range("C4").Value = [double]"1.233997095605541e-18"
$ws.Range("D4").Value = [double]"2.176765515339244e-19"
$ws.Range("A5").Value = [double]"0.02773906721555619"
$ws.Range("B5").Value = [double]"0.1955725434464861"
$ws.Range("C5").Value = [double]"0.2437218137074627"
$ws.Range("D5").Value = [double]"0.7366850471976923"

$ws = $wb.Worksheets.Item(3)
$ws.Range("A1").Value = [double]"0.8379594278711199"
$ws.Range("B1").Value = [double]"0"
$ws.Range("A2").Value = [double]"0.1459159178751359"
$ws.Range("B2").Value = [double]"0.8793706712783951"
$ws.Range("C2").Value = [double]"0.5966131359133103"
$ws.Range("D2").Value = [double]"0.0002325062982800358"
$ws.Range("A3").Value = [double]"0.01260174915820445"
$ws.Range("B3").Value = [double]"0.09424433784111469"
$ws.Range("C3").Value = [double]"0.3142743015993417"
$ws.Range("D3").Value = [double]"0.6410733375532822"
$ws.Range("A4").Value = [double]"0.0005592594554913408"
$ws.Range("B4").Value = [double]"0.004180557277934698"
$ws.Range("C4").Value = [double]"0.01388410901487574"
$ws.Range("D4").Value = [double]"0.02475676119273205"
$ws.Range("A5").Value = [double]"0.002963630738887504"
$ws.Range("B5").Value = [double]"0.0222044559542977"
$ws.Range("C5").Value = [double]"0.0752284385713113"
$ws.Range("D5").Value = [double]"0.3339373920122664"

$ws = $wb.Worksheets.Item(4)
$ws.Range("A1").Value = [double]"0.423359647032813"
$ws.Range("A2").Value = [double]"0.5742192490126857"
$ws.Range("B2").Value = [double]"0.9945783579600974"
$ws.Range("C2").Value = [double]"0.4914830316018843"
$ws.Range("A3").Value = [double]"0"
$ws.Range("B3").Value = [double]"0"
$ws.Range("C3").Value = [double]"0"
$ws.Range("D3").Value = [double]"0"
$ws.Range("A4").Value = [double]"1.35686056107676e-11"
$ws.Range("B4").Value = [double]"2.69568457256313e-11"
$ws.Range("C4").Value = [double]"1.39540747693782e-10"
$ws.Range("D4").Value = [double]"0"
$ws.Range("A5").Value = [double]"0.002421126292675022"
$ws.Range("B5").Value = [double]"0.005421664364687924"
$ws.Range("C5").Value = [double]"0.508516953357414"
$ws.Range("D5").Value = [double]"0.9999999719507564"

$ws = $wb.Worksheets.Item(5)
$ws.Range("A1").Value = [double]"0.423359647032813"
$ws.Range("A2").Value = [double]"0.5742192490126857"
$ws.Range("B2").Value = [double]"0.9945783579600974"
$ws.Range("C2").Value = [double]"0.4914830316018843"
$ws.Range("A3").Value = [double]"0"
$ws.Range("B3").Value = [double]"0"
$ws.Range("C3").Value = [double]"0"
$ws.Range("D3").Value = [double]"0"
$ws.Range("A4").Value = [double]"1.35686056107676e-11"
$ws.Range("B4").Value = [double]"2.69568457256313e-11"
$ws.Range("C4").Value = [double]"1.39540747693782e-10"
$ws.Range("D4").Value = [double]"0"
$ws.Range("A5").Value = [double]"0.002421126292675022"
$ws.Range("B5").Value = [double]"0.005421664364687924"
$ws.Range("C5").Value = [double]"0.508516953357414"
$ws.Range("D5").Value = [double]"0.9999999719507564"

$ws = $wb.Worksheets.Item(6)
$ws.Range("A1").Value = [double]"0.6993333887896118"
$ws.Range("B1").Value = [double]"0.1267966844611701"
$ws.Range("A2").Value = [double]"0.2708263575960009"
$ws.Range("B2").Value = [double]"0.6816187700345274"
$ws.Range("C2").Value = [double]"0.04932949192737977"
$ws.Range("D2").Value = [double]"0"
$ws.Range("A3").Value = [double]"0"
$ws.Range("B3").Value = [double]"0"
$ws.Range("C3").Value = [double]"0"
$ws.Range("D3").Value = [double]"0"
$ws.Range("A4").Value = [double]"0.02980577981898093"
$ws.Range("B4").Value = [double]"0.191361135335625"
$ws.Range("C4").Value = [double]"0.9494964142238367"
$ws.Range("D4").Value = [double]"0.4775888994233586"
$ws.Range("A5").Value = [double]"3.445889424544195e-05"
$ws.Range("B5").Value = [double]"0.0002233952675162994"
$ws.Range("C5").Value = [double]"0.001174078947622709"
$ws.Range("D5").Value = [double]"0.5224110725273976"

$ws = $wb.Worksheets.Item(7)
$ws.Range("A1").Value = [double]"0.6993333887896118"
$ws.Range("B1").Value = [double]"0.1267966844611701"
$ws.Range("A2").Value = [double]"0.2708263575960009"
$ws.Range("B2").Value = [double]"0.6816187700345274"
$ws.Range("C2").Value = [double]"0.04932949192737977"
$ws.Range("D2").Value = [double]"0"
$ws.Range("A3").Value = [double]"0"
$ws.Range("B3").Value = [double]"0"
$ws.Range("C3").Value = [double]"0"
$ws.Range("D3").Value = [double]"0"
$ws.Range("A4").Value = [double]"0.02980577981898093"
$ws.Range("B4").Value = [double]"0.191361135335625"
$ws.Range("C4").Value = [double]"0.9494964142238367"
$ws.Range("D4").Value = [double]"0.4775888994233586"
$ws.Range("A5").Value = [double]"3.445889424544195e-05"
$ws.Range("B5").Value = [double]"0.0002233952675162994"
$ws.Range("C5").Value = [double]"0.001174078947622709"
$ws.Range("D5").Value = [double]"0.5224110725273976"

$ws = $wb.Worksheets.Item(8)
$ws.Range("A1").Value = [double]"0"
$ws.Range("B1").Value = [double]"0"
$ws.Range("A2").Value = [double]"0.9849787172158034"
$ws.Range("B2").Value = [double]"0.474413655080572"
$ws.Range("C2").Value = [double]"0"
$ws.Range("A3").Value = [double]"1.071793831925002e-21"
$ws.Range("B3").Value = [double]"3.786611284348519e-20"
$ws.Range("C3").Value = [double]"7.971976402181327e-20"
$ws.Range("D3").Value = [double]"9.300470836452734e-19"
$ws.Range("A4").Value = [double]"0.01451627633583374"
$ws.Range("B4").Value = [double]"0.5077451870536952"
$ws.Range("C4").Value = [double]"0.962438836884692"
$ws.Range("D4").Value = [double]"0.5617945309016854"
$ws.Range("A5").Value = [double]"0.0005049915472018728"
$ws.Range("B5").Value = [double]"0.0178411802174752"
$ws.Range("C5").Value = [double]"0.03756114821414708"
$ws.Range("D5").Value = [double]"0.4382054661548754"
